$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados..." timestamp shown in A1 (last shared string in the sst)
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 01:17"

# Refreshed COVID-19 country figures (columns: Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
# The data refresh also nudges a handful of countries past their neighbours in the
# "Casos totales" ranking, so a few rows swap which country they display.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 4858654
$ws.Cells.Item(4, 3).Value = 45626
$ws.Cells.Item(4, 4).Value = 2442684
$ws.Cells.Item(4, 5).Value = 2257077
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 528
$ws.Cells.Item(4, 8).Value = 158893

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 2751665
$ws.Cells.Item(5, 3).Value = 17988
$ws.Cells.Item(5, 4).Value = 1912319
$ws.Cells.Item(5, 5).Value = 744644
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 572
$ws.Cells.Item(5, 8).Value = 94702

# Row 10: Peru
$ws.Cells.Item(10, 1).Value = "Peru"
$ws.Cells.Item(10, 2).Value = 433100
$ws.Cells.Item(10, 3).Value = 4250
$ws.Cells.Item(10, 4).Value = 298091
$ws.Cells.Item(10, 5).Value = 115198
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 197
$ws.Cells.Item(10, 8).Value = 19811

# Row 13: Colombia
$ws.Cells.Item(13, 1).Value = "Colombia"
$ws.Cells.Item(13, 2).Value = 327850
$ws.Cells.Item(13, 3).Value = 10199
$ws.Cells.Item(13, 4).Value = 173727
$ws.Cells.Item(13, 5).Value = 143106
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 367
$ws.Cells.Item(13, 8).Value = 11017

# Row 21: Alemania
$ws.Cells.Item(21, 1).Value = "Alemania"
$ws.Cells.Item(21, 2).Value = 212320
$ws.Cells.Item(21, 3).Value = 858
$ws.Cells.Item(21, 4).Value = 193600
$ws.Cells.Item(21, 5).Value = 9488
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 6
$ws.Cells.Item(21, 8).Value = 9232

# Row 25: Canada
$ws.Cells.Item(25, 1).Value = "Canada"
$ws.Cells.Item(25, 2).Value = 117031
$ws.Cells.Item(25, 3).Value = 147
$ws.Cells.Item(25, 4).Value = 101597
$ws.Cells.Item(25, 5).Value = 6487
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = 8947

# Row 29: Egipto
$ws.Cells.Item(29, 1).Value = "Egipto"
$ws.Cells.Item(29, 2).Value = 94640
$ws.Cells.Item(29, 3).Value = 157
$ws.Cells.Item(29, 4).Value = 44066
$ws.Cells.Item(29, 5).Value = 45686
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 23
$ws.Cells.Item(29, 8).Value = 4888

# Row 50: Nigeria
$ws.Cells.Item(50, 1).Value = "Nigeria"
$ws.Cells.Item(50, 2).Value = 44129
$ws.Cells.Item(50, 3).Value = 288
$ws.Cells.Item(50, 4).Value = 20663
$ws.Cells.Item(50, 5).Value = 22570
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 8
$ws.Cells.Item(50, 8).Value = 896

# Row 52: Barein
$ws.Cells.Item(52, 1).Value = "Barein"
$ws.Cells.Item(52, 2).Value = 41835
$ws.Cells.Item(52, 3).Value = 299
$ws.Cells.Item(52, 4).Value = 39007
$ws.Cells.Item(52, 5).Value = 2678
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 3
$ws.Cells.Item(52, 8).Value = 150

# Row 54: Japon
$ws.Cells.Item(54, 1).Value = "Japon"
$ws.Cells.Item(54, 2).Value = 38687
$ws.Cells.Item(54, 3).Value = 1998
$ws.Cells.Item(54, 4).Value = 26487
$ws.Cells.Item(54, 5).Value = 11188
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 1012

# Row 55: Ghana
$ws.Cells.Item(55, 1).Value = "Ghana"
$ws.Cells.Item(55, 2).Value = 37812
$ws.Cells.Item(55, 3).Value = 798
$ws.Cells.Item(55, 4).Value = 34313
$ws.Cells.Item(55, 5).Value = 3308
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 9
$ws.Cells.Item(55, 8).Value = 191

# Row 56: Kirguistan
$ws.Cells.Item(56, 1).Value = "Kirguistan"
$ws.Cells.Item(56, 2).Value = 37129
$ws.Cells.Item(56, 3).Value = 410
$ws.Cells.Item(56, 4).Value = 27927
$ws.Cells.Item(56, 5).Value = 7782
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 11
$ws.Cells.Item(56, 8).Value = 1420

# Row 57: Afganistan
$ws.Cells.Item(57, 1).Value = "Afganistan"
$ws.Cells.Item(57, 2).Value = 36747
$ws.Cells.Item(57, 3).Value = 37
$ws.Cells.Item(57, 4).Value = 25669
$ws.Cells.Item(57, 5).Value = 9790
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 4
$ws.Cells.Item(57, 8).Value = 1288

# Row 81: Bulgaria
$ws.Cells.Item(81, 1).Value = "Bulgaria"
$ws.Cells.Item(81, 2).Value = 12159
$ws.Cells.Item(81, 3).Value = 204
$ws.Cells.Item(81, 4).Value = 6684
$ws.Cells.Item(81, 5).Value = 5071
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 16
$ws.Cells.Item(81, 8).Value = 404

# Row 86: Noruega
$ws.Cells.Item(86, 1).Value = "Noruega"
$ws.Cells.Item(86, 2).Value = 9334
$ws.Cells.Item(86, 3).Value = 66
$ws.Cells.Item(86, 4).Value = 8752
$ws.Cells.Item(86, 5).Value = 326
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 256

# Row 98: Paraguay
$ws.Cells.Item(98, 1).Value = "Paraguay"
$ws.Cells.Item(98, 2).Value = 5724
$ws.Cells.Item(98, 3).Value = 80
$ws.Cells.Item(98, 4).Value = 4249
$ws.Cells.Item(98, 5).Value = 1420
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 3
$ws.Cells.Item(98, 8).Value = 55

# Row 109: Zimbabue
$ws.Cells.Item(109, 1).Value = "Zimbabue"
$ws.Cells.Item(109, 2).Value = 4075
$ws.Cells.Item(109, 3).Value = 154
$ws.Cells.Item(109, 4).Value = 1057
$ws.Cells.Item(109, 5).Value = 2938
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 10
$ws.Cells.Item(109, 8).Value = 80

# Row 128: Ruanda
$ws.Cells.Item(128, 1).Value = "Ruanda"
$ws.Cells.Item(128, 2).Value = 2092
$ws.Cells.Item(128, 3).Value = 30
$ws.Cells.Item(128, 4).Value = 1169
$ws.Cells.Item(128, 5).Value = 918
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 5

# Row 129: Estonia
$ws.Cells.Item(129, 1).Value = "Estonia"
$ws.Cells.Item(129, 2).Value = 2080
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 4).Value = 1935
$ws.Cells.Item(129, 5).Value = 82
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 63

# Row 139: Uruguay
$ws.Cells.Item(139, 1).Value = "Uruguay"
$ws.Cells.Item(139, 2).Value = 1291
$ws.Cells.Item(139, 3).Value = 5
$ws.Cells.Item(139, 4).Value = 1023
$ws.Cells.Item(139, 5).Value = 232
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 36

# Row 147: Niger
$ws.Cells.Item(147, 1).Value = "Niger"
$ws.Cells.Item(147, 2).Value = 1152
$ws.Cells.Item(147, 3).Value = 5
$ws.Cells.Item(147, 4).Value = 1032
$ws.Cells.Item(147, 5).Value = 51
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 69

# Row 148: Burkina Faso
$ws.Cells.Item(148, 1).Value = "Burkina Faso"
$ws.Cells.Item(148, 2).Value = 1150
$ws.Cells.Item(148, 3).Value = 7
$ws.Cells.Item(148, 4).Value = 947
$ws.Cells.Item(148, 5).Value = 150
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 53

# Row 149: Togo
$ws.Cells.Item(149, 1).Value = "Togo"
$ws.Cells.Item(149, 2).Value = 976
$ws.Cells.Item(149, 3).Value = 15
$ws.Cells.Item(149, 4).Value = 663
$ws.Cells.Item(149, 5).Value = 294
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 19

# Row 160: Bahamas
$ws.Cells.Item(160, 1).Value = "Bahamas"
$ws.Cells.Item(160, 2).Value = 679
$ws.Cells.Item(160, 3).Value = 31
$ws.Cells.Item(160, 4).Value = 91
$ws.Cells.Item(160, 5).Value = 574
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 14

# Row 161: Reunion
$ws.Cells.Item(161, 1).Value = "Reunion"
$ws.Cells.Item(161, 2).Value = 667
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 592
$ws.Cells.Item(161, 5).Value = 71
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 4

# Row 162: Vietnam
$ws.Cells.Item(162, 1).Value = "Vietnam"
$ws.Cells.Item(162, 2).Value = 652
$ws.Cells.Item(162, 3).Value = 32
$ws.Cells.Item(162, 4).Value = 374
$ws.Cells.Item(162, 5).Value = 272
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 6

# Row 173: Eritrea
$ws.Cells.Item(173, 1).Value = "Eritrea"
$ws.Cells.Item(173, 2).Value = 282
$ws.Cells.Item(173, 3).Value = 3
$ws.Cells.Item(173, 4).Value = 225
$ws.Cells.Item(173, 5).Value = 57
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0
